# Temporary change pending resolution of NCIOCPL/cgov-digital-platform#1541:
# Change language toggles to look for "Spanish" instead of "Espanol".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pages_with_translation")

# Column E ("expectedtoggle") for data rows 2-8 previously expected "Espanol";
# update them to expect "Spanish" instead.
$ws.Range("E2:E8").Value = "Spanish"

# Matches the author's final selection left in the sheet.
$ws.Range("E5").Select() | Out-Null
